$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.88
$ws.Range("I2").Value = 2.7
$ws.Range("J2").Value = 3.75
$ws.Range("L2").Value = 3.5
$ws.Range("W2").Value = 7
$ws.Range("X2").Value = 12
$ws.Range("Y2").Value = 12
$ws.Range("Z2").Value = 29
$ws.Range("AA2").Value = 29
$ws.Range("AH2").Value = 6.5
$ws.Range("AI2").Value = 11
$ws.Range("AJ2").Value = 11
$ws.Range("AK2").Value = 26
$ws.Range("AL2").Value = 26
$ws.Range("AN2").Value = 4.75
$ws.Range("AW2").Value = 4.5
# Row 3
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
# Row 4
$ws.Range("I4").Value = 2.05
$ws.Range("AJ4").Value = 9.5
# Row 5
$ws.Range("G5").Value = 1.29
$ws.Range("H5").Value = 5
$ws.Range("K5").Value = 2.5
$ws.Range("L5").Value = 8.5
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 9.5
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 3.5
$ws.Range("Q5").Value = 1.93
$ws.Range("R5").Value = 1.93
$ws.Range("S5").Value = 1.33
$ws.Range("T5").Value = 3.25
$ws.Range("U5").Value = 2.5
$ws.Range("V5").Value = 1.5
$ws.Range("AC5").Value = 9.5
$ws.Range("AN5").Value = 3.2
$ws.Range("AP5").Value = 19
$ws.Range("AR5").Value = 41
$ws.Range("AS5").Value = 151
$ws.Range("AT5").Value = 3.25
$ws.Range("AU5").Value = 10
$ws.Range("AV5").Value = 67
$ws.Range("AX5").Value = 41
$ws.Range("AY5").Value = 41
$ws.Range("AZ5").Value = 201
$ws.Range("BA5").Value = 201
# Row 6
$ws.Range("G6").Value = 1.45
$ws.Range("K6").Value = 2.25
$ws.Range("N6").Value = 9.5
$ws.Range("Q6").Value = 2.03
$ws.Range("R6").Value = 1.83
$ws.Range("Y6").Value = 9
$ws.Range("AC6").Value = 9.5
$ws.Range("AD6").Value = 8
$ws.Range("AF6").Value = 81
$ws.Range("AH6").Value = 15
$ws.Range("AI6").Value = 34
$ws.Range("AJ6").Value = 21
$ws.Range("AS6").Value = 151
# Row 7
$ws.Range("G7").Value = 1.3
$ws.Range("H7").Value = 4.5
$ws.Range("I7").Value = 12
$ws.Range("J7").Value = 1.83
$ws.Range("K7").Value = 2.25
$ws.Range("L7").Value = 11
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 10
$ws.Range("O7").Value = 1.33
$ws.Range("P7").Value = 3.25
$ws.Range("Q7").Value = 2.08
$ws.Range("R7").Value = 1.73
$ws.Range("S7").Value = 1.44
$ws.Range("T7").Value = 2.63
$ws.Range("U7").Value = 2.75
$ws.Range("V7").Value = 1.4
$ws.Range("X7").Value = 5
$ws.Range("Z7").Value = 7.5
$ws.Range("AC7").Value = 8
$ws.Range("AD7").Value = 9.5
$ws.Range("AE7").Value = 34
$ws.Range("AF7").Value = 126
$ws.Range("AH7").Value = 21
$ws.Range("AI7").Value = 51
$ws.Range("AJ7").Value = 34
$ws.Range("AK7").Value = 151
$ws.Range("AL7").Value = 101
$ws.Range("AM7").Value = 101
$ws.Range("AN7").Value = 3
$ws.Range("AO7").Value = 6.5
$ws.Range("AQ7").Value = 19
$ws.Range("AT7").Value = 2.63
$ws.Range("AU7").Value = 12
$ws.Range("AV7").Value = 101
$ws.Range("AW7").Value = 11
$ws.Range("AX7").Value = 51
$ws.Range("AZ7").Value = 351
# Row 8
$ws.Range("H8").Value = 4.5
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("U8").Value = 2.1
$ws.Range("V8").Value = 1.67
$ws.Range("W8").Value = 6.5
$ws.Range("Z8").Value = 8.5
$ws.Range("AH8").Value = 21
$ws.Range("AJ8").Value = 26
$ws.Range("AW8").Value = 9
# Row 9
$ws.Range("Q9").Value = 1.9
$ws.Range("R9").Value = 1.95
# Row 11
$ws.Range("G11").Value = 1.65
$ws.Range("H11").Value = 3.6
$ws.Range("I11").Value = 5.75
$ws.Range("J11").Value = 2.3
$ws.Range("K11").Value = 2.1
$ws.Range("L11").Value = 5.5
$ws.Range("U11").Value = 2
$ws.Range("V11").Value = 1.73
$ws.Range("X11").Value = 7
$ws.Range("Y11").Value = 8.5
$ws.Range("Z11").Value = 12
$ws.Range("AB11").Value = 29
$ws.Range("AC11").Value = 8.5
$ws.Range("AD11").Value = 7
$ws.Range("AE11").Value = 17
$ws.Range("AF11").Value = 51
$ws.Range("AH11").Value = 13
$ws.Range("AI11").Value = 29
$ws.Range("AJ11").Value = 19
$ws.Range("AK11").Value = 51
$ws.Range("AL11").Value = 41
$ws.Range("AN11").Value = 3.5
$ws.Range("AO11").Value = 8.5
$ws.Range("AQ11").Value = 29
$ws.Range("AS11").Value = 151
$ws.Range("AU11").Value = 9
$ws.Range("AW11").Value = 7
$ws.Range("AX11").Value = 29
$ws.Range("AZ11").Value = 101
$ws.Range("BA11").Value = 151
$ws.Range("BB11").Value = 301

Write-Host "Updated cells"
